# Update "想去人数" (interested-people count) figures on the 展览 (Exhibition)
# and 全部类型 (All types) sheets, reflecting a refreshed data pull.

$wb = $excel.ActiveWorkbook

$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F4").Value = 3
$wsExhibition.Range("F5").Value = 3788
$wsExhibition.Range("F8").Value = 214
$wsExhibition.Range("F9").Value = 4

$wsAllTypes = $wb.Worksheets.Item("全部类型")
$wsAllTypes.Range("F8").Value = 3
$wsAllTypes.Range("F9").Value = 3788
$wsAllTypes.Range("F13").Value = 214
$wsAllTypes.Range("F14").Value = 4
